$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 174.75
$ws.Range("I4").Value = 133.33333
$ws.Range("J4").Value = 299
$ws.Range("K4").Value = 133.33333
$ws.Range("L4").Value = 299
$ws.Range("M4").Value = -19.33332999999999
$ws.Range("N4").Value = -527

# Row 32
$ws.Range("H32").Value = 700
$ws.Range("I32").Value = 433.33334
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 433.33334
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = -107.33334
$ws.Range("N32").Value = -2152

# Row 103
$ws.Range("H103").Value = 227381.19
$ws.Range("I103").Value = 500079
$ws.Range("K103").Value = 1500237
$ws.Range("M103").Value = -1499651

# Row 112
$ws.Range("H112").Value = 3585283.8
$ws.Range("J112").Value = 3585283.8
$ws.Range("L112").Value = 10755851.4
$ws.Range("N112").Value = -10758067.4

# Row 132
$ws.Range("H132").Value = 3516.2693
$ws.Range("I132").Value = 4025.0952
$ws.Range("K132").Value = 12075.2856
$ws.Range("M132").Value = -9545.285600000001

# Row 137
$ws.Range("H137").Value = 1880.25
$ws.Range("I137").Value = 1910.6
$ws.Range("J137").Value = 1771.8572
$ws.Range("K137").Value = 5731.799999999999
$ws.Range("L137").Value = 5315.571599999999
$ws.Range("M137").Value = -3181.799999999999
$ws.Range("N137").Value = -10415.5716

# Row 138
$ws.Range("H138").Value = 19609832
$ws.Range("I138").Value = 27028082
$ws.Range("J138").Value = 4455.5
$ws.Range("K138").Value = 81084246
$ws.Range("L138").Value = 13366.5
$ws.Range("M138").Value = -81079106
$ws.Range("N138").Value = -23646.5

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 18999.5
$ws.Range("J24").Value = 18999.5
$ws.Range("L24").Value = 18999.5
$ws.Range("N24").Value = -19747.5

# Row 55
$ws.Range("H55").Value = 33333.332
$ws.Range("J55").Value = 33333.332
$ws.Range("L55").Value = 33333.332
$ws.Range("N55").Value = -33963.332

# Row 100
$ws.Range("H100").Value = 18999.5
$ws.Range("J100").Value = 18999.5
$ws.Range("L100").Value = 18999.5
$ws.Range("N100").Value = -21163.5

$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 31821.5
$ws.Range("J100").Value = 31821.5
$ws.Range("L100").Value = 31821.5
$ws.Range("N100").Value = -33985.5

# Row 134
$ws.Range("H134").Value = 3288.9722
$ws.Range("I134").Value = 3633.0967
$ws.Range("J134").Value = 1155.4
$ws.Range("K134").Value = 10899.2901
$ws.Range("L134").Value = 3466.2
$ws.Range("M134").Value = -8364.2901
$ws.Range("N134").Value = -8536.200000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4061.718
$ws.Range("I31").Value = 2624.1667
$ws.Range("J31").Value = 6361.8
$ws.Range("K31").Value = 2624.1667
$ws.Range("L31").Value = 6361.8
$ws.Range("M31").Value = -2329.1667
$ws.Range("N31").Value = -6951.8

# Row 34
$ws.Range("H34").Value = 4061.718
$ws.Range("I34").Value = 2624.1667
$ws.Range("J34").Value = 6361.8
$ws.Range("K34").Value = 2624.1667
$ws.Range("L34").Value = 6361.8
$ws.Range("M34").Value = -2422.1667
$ws.Range("N34").Value = -6765.8

# Row 50
$ws.Range("H50").Value = 19990
$ws.Range("J50").Value = 19990
$ws.Range("L50").Value = 19990
$ws.Range("N50").Value = -21240

# Row 86
$ws.Range("H86").Value = 9271518
$ws.Range("I86").Value = 5814.154
$ws.Range("J86").Value = 33362348
$ws.Range("K86").Value = 5814.154
$ws.Range("L86").Value = 33362348
$ws.Range("M86").Value = -4691.154
$ws.Range("N86").Value = -33364594

# Row 89
$ws.Range("H89").Value = 9271518
$ws.Range("I89").Value = 5814.154
$ws.Range("J89").Value = 33362348
$ws.Range("K89").Value = 29070.77
$ws.Range("L89").Value = 166811740
$ws.Range("M89").Value = -23454.77
$ws.Range("N89").Value = -166822972

# Row 132
$ws.Range("H132").Value = 2929.1
$ws.Range("I132").Value = 1944.1154
$ws.Range("K132").Value = 5832.3462
$ws.Range("M132").Value = -3302.3462

# Row 134
$ws.Range("H134").Value = 1056.5
$ws.Range("I134").Value = 970.7273
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2912.1819
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -377.1819
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4210
$ws.Range("I3").Value = 1923
$ws.Range("K3").Value = 5769
$ws.Range("M3").Value = -5657

# Row 5
$ws.Range("H5").Value = 1010.6316
$ws.Range("I5").Value = 395.25
$ws.Range("J5").Value = 2065.5715
$ws.Range("K5").Value = 1185.75
$ws.Range("L5").Value = 6196.7145
$ws.Range("M5").Value = -1073.75
$ws.Range("N5").Value = -6420.7145

# Row 116
$ws.Range("H116").Value = 1884.2858
$ws.Range("J116").Value = 1896.6666
$ws.Range("L116").Value = 5689.9998
$ws.Range("N116").Value = -12573.9998

# Row 129
$ws.Range("H129").Value = 1362.68
$ws.Range("I129").Value = 961.8182
$ws.Range("K129").Value = 2885.4546
$ws.Range("M129").Value = 2114.5454

# Row 131
$ws.Range("H131").Value = 784.13
$ws.Range("J131").Value = 796.97894
$ws.Range("L131").Value = 2390.93682
$ws.Range("N131").Value = -12470.93682

# Row 133
$ws.Range("H133").Value = 6000
$ws.Range("J133").Value = 6000
$ws.Range("L133").Value = 18000
$ws.Range("N133").Value = -28120

# Row 135
$ws.Range("H135").Value = 1010.6316
$ws.Range("I135").Value = 395.25
$ws.Range("J135").Value = 2065.5715
$ws.Range("K135").Value = 3557.25
$ws.Range("L135").Value = 18590.1435
$ws.Range("M135").Value = -1022.25
$ws.Range("N135").Value = -23660.1435

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 13000
$ws.Range("J5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("N5").Value = -13224

# Row 80
$ws.Range("H80").Value = 3334.9119
$ws.Range("I80").Value = 2465.1538
$ws.Range("J80").Value = 3873.3333
$ws.Range("K80").Value = 2465.1538
$ws.Range("L80").Value = 3873.3333
$ws.Range("M80").Value = -1467.1538
$ws.Range("N80").Value = -5869.3333

# Row 83
$ws.Range("H83").Value = 3334.9119
$ws.Range("I83").Value = 2465.1538
$ws.Range("J83").Value = 3873.3333
$ws.Range("K83").Value = 12325.769
$ws.Range("L83").Value = 19366.6665
$ws.Range("M83").Value = -7333.769
$ws.Range("N83").Value = -29350.6665

# Row 132
$ws.Range("H132").Value = 82715.71000000001
$ws.Range("I132").Value = 17253
$ws.Range("K132").Value = 51759
$ws.Range("M132").Value = -49229

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 16131015
$ws.Range("I136").Value = 25000984
$ws.Range("J136").Value = 3800.2273
$ws.Range("K136").Value = 75002952
$ws.Range("L136").Value = 11400.6819
$ws.Range("M136").Value = -75000402
$ws.Range("N136").Value = -16500.6819
